$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from G1 into H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save column values per row (0/1), matching the diff row by row
$values = @(0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
